$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the gender value in C2 from "M" to lowercase "m"
$ws.Range("C2").Value = "m"

# Move the active selection to E5 (no data, just the final cursor position)
$ws.Range("E5").Select()
